# BookingDetails sheet: update the stored Password test-data value
# (B2) from the old generated password to the new one, and tidy up
# the two trailing helper columns (S:T) that accompanied it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Password value used by the Cucumber automation data file.
$ws.Range("B2").Value = "E9Z51L"

# S2/T2 were stray blank, bordered placeholder cells - remove them.
$ws.Range("S2:T2").Clear()

# S1/T1 remain as blank cells but pick up the plain centered/no-border
# "section heading" look (same font/alignment as the other un-bordered
# header cells) instead of their old bordered style.
$ws.Range("S1:T1").Font.Name = "Calibri"
$ws.Range("S1:T1").Font.Size = 11
$ws.Range("S1:T1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("S1:T1").HorizontalAlignment = -4108
$ws.Range("S1:T1").VerticalAlignment = -4108

# Leave the selection on the cell that was actually edited.
[void]$ws.Range("B2").Select()
